$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.036334996200109
$ws.Range("D2").Value = 1.04053603843413
$ws.Range("E2").Value = 1.04409588715953
$ws.Range("F2").Value = 1.052215392369115
$ws.Range("I2").Value = 1.03971543168033
$ws.Range("J2").Value = 1.041443875618211
$ws.Range("K2").Value = 1.043318109209606
$ws.Range("L2").Value = 1.046867900200587
$ws.Range("M2").Value = 1.054964741209956
$ws.Range("N2").Value = 1.042922845376871
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037315487867556
$ws.Range("D3").Value = 1.041106513308931
$ws.Range("E3").Value = 1.045002393081642
$ws.Range("F3").Value = 1.053287367826528
$ws.Range("I3").Value = 1.039935749525724
$ws.Range("J3").Value = 1.042068132521718
$ws.Range("K3").Value = 1.043699658927979
$ws.Range("L3").Value = 1.047585325712558
$ws.Range("M3").Value = 1.055848852529881
$ws.Range("N3").Value = 1.043547988796784
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037950091557073
$ws.Range("D4").Value = 1.04147531272245
$ws.Range("E4").Value = 1.045589542672901
$ws.Range("F4").Value = 1.053981910279077
$ws.Range("I4").Value = 1.040076780118256
$ws.Range("J4").Value = 1.04247162707525
$ws.Range("K4").Value = 1.043945509319753
$ws.Range("L4").Value = 1.048049495865351
$ws.Range("M4").Value = 1.056421235249247
$ws.Range("N4").Value = 1.043952056358864
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038216916730792
$ws.Range("D5").Value = 1.041630273961693
$ws.Range("E5").Value = 1.04583651822982
$ws.Range("F5").Value = 1.054274111257837
$ws.Range("I5").Value = 1.040135702486746
$ws.Range("J5").Value = 1.042641149839911
$ws.Range("K5").Value = 1.044048615253205
$ws.Range("L5").Value = 1.048244619806119
$ws.Range("M5").Value = 1.056661936891508
$ws.Range("N5").Value = 1.044121819865295
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038261720051903
$ws.Range("D6").Value = 1.041656287778405
$ws.Range("E6").Value = 1.045877994571123
$ws.Range("F6").Value = 1.054323185692882
$ws.Range("I6").Value = 1.040145574277725
$ws.Range("J6").Value = 1.042669607202317
$ws.Range("K6").Value = 1.044065912525405
$ws.Range("L6").Value = 1.048277381180097
$ws.Range("M6").Value = 1.056702355954195
$ws.Range("N6").Value = 1.04415031764042
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037953656740297
$ws.Range("D7").Value = 1.041477383646479
$ws.Range("E7").Value = 1.045592842231779
$ws.Range("F7").Value = 1.053985813837706
$ws.Range("I7").Value = 1.040077568884295
$ws.Range("J7").Value = 1.04247389266544
$ws.Range("K7").Value = 1.043946888008347
$ws.Range("L7").Value = 1.048052103173726
$ws.Range("M7").Value = 1.056424451234584
$ws.Range("N7").Value = 1.043954325166452
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036666324500924
$ws.Range("D8").Value = 1.040728901695388
$ws.Range("E8").Value = 1.044402124521028
$ws.Range("F8").Value = 1.052577484134245
$ws.Range("I8").Value = 1.039790205484709
$ws.Range("J8").Value = 1.041654937153917
$ws.Range("K8").Value = 1.04344726993143
$ws.Range("L8").Value = 1.047110367876507
$ws.Range("M8").Value = 1.0552634674737
$ws.Range("N8").Value = 1.043134206644167
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034399129034715
$ws.Range("D9").Value = 1.039407478026022
$ws.Range("E9").Value = 1.042308406590031
$ws.Range("F9").Value = 1.050102773980514
$ws.Range("I9").Value = 1.039272147007576
$ws.Range("J9").Value = 1.040208481669127
$ws.Range("K9").Value = 1.042558980092461
$ws.Range("L9").Value = 1.045450541934604
$ws.Range("M9").Value = 1.053220022745829
$ws.Range("N9").Value = 1.041685697026677
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032888533276178
$ws.Range("D10").Value = 1.038524940308827
$ws.Range("E10").Value = 1.040915657467852
$ws.Range("F10").Value = 1.04845768091539
$ws.Range("I10").Value = 1.038918954890246
$ws.Range("J10").Value = 1.039241963007085
$ws.Range("K10").Value = 1.041961541759487
$ws.Range("L10").Value = 1.044343784427034
$ws.Range("M10").Value = 1.051859352266882
$ws.Range("N10").Value = 1.040717805797275
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032234639743699
$ws.Range("D11").Value = 1.038142434594602
$ws.Range("E11").Value = 1.040313318529861
$ws.Range("F11").Value = 1.047746463921254
$ws.Range("I11").Value = 1.038764171849902
$ws.Range("J11").Value = 1.038822931959842
$ws.Range("K11").Value = 1.041701613917158
$ws.Range("L11").Value = 1.043864505204257
$ws.Range("M11").Value = 1.051270560723082
$ws.Range("N11").Value = 1.040298179677881
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031991785611178
$ws.Range("D12").Value = 1.038000302257964
$ws.Range("E12").Value = 1.040089693572564
$ws.Range("F12").Value = 1.047482455011003
$ws.Range("I12").Value = 1.038706401400363
$ws.Range("J12").Value = 1.038667207281624
$ws.Range("K12").Value = 1.041604880945703
$ws.Range("L12").Value = 1.043686473251993
$ws.Range("M12").Value = 1.051051916169469
$ws.Range("N12").Value = 1.040142233852758
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032043877237026
$ws.Range("D13").Value = 1.038030792485696
$ws.Range("E13").Value = 1.040137656870249
$ws.Range("F13").Value = 1.047539078172472
$ws.Range("I13").Value = 1.038718805889575
$ws.Range("J13").Value = 1.038700614269616
$ws.Range("K13").Value = 1.041625638814778
$ws.Range("L13").Value = 1.043724661986129
$ws.Range("M13").Value = 1.051098813496296
$ws.Range("N13").Value = 1.040175688282504
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032214564710214
$ws.Range("D14").Value = 1.038130686957039
$ws.Range("E14").Value = 1.040294831361594
$ws.Range("F14").Value = 1.047724637416442
$ws.Range("I14").Value = 1.03875940217666
$ws.Range("J14").Value = 1.038810061295126
$ws.Range("K14").Value = 1.041693621691012
$ws.Range("L14").Value = 1.043849789152837
$ws.Range("M14").Value = 1.051252486292327
$ws.Range("N14").Value = 1.040285290735344
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032319735018899
$ws.Range("D15").Value = 1.038192228303286
$ws.Range("E15").Value = 1.04039168642439
$ws.Range("F15").Value = 1.047838988969859
$ws.Range("I15").Value = 1.038784378183088
$ws.Range("J15").Value = 1.038877484914959
$ws.Range("K15").Value = 1.041735483808494
$ws.Range("L15").Value = 1.043926883311594
$ws.Range("M15").Value = 1.051347176988449
$ws.Range("N15").Value = 1.040352810104451
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032931934371084
$ws.Range("D16").Value = 1.03855031847973
$ws.Range("E16").Value = 1.040955648179013
$ws.Range("F16").Value = 1.048504905676159
$ws.Range("I16").Value = 1.038929188434432
$ws.Range("J16").Value = 1.039269761764635
$ws.Range("K16").Value = 1.041978766425075
$ws.Range("L16").Value = 1.044375591693477
$ws.Range("M16").Value = 1.051898436642035
$ws.Range("N16").Value = 1.040745644032249
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033316005807669
$ws.Range("D17").Value = 1.038774843426238
$ws.Range("E17").Value = 1.04130960264434
$ws.Range("F17").Value = 1.04892291769189
$ws.Range("I17").Value = 1.039019529589542
$ws.Range("J17").Value = 1.03951568727794
$ws.Range("K17").Value = 1.04213104179818
$ws.Range("L17").Value = 1.044657042729569
$ws.Range("M17").Value = 1.05224433148486
$ws.Range("N17").Value = 1.04099191878799
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033540047877939
$ws.Range("D18").Value = 1.038905770087671
$ws.Range("E18").Value = 1.041516128738066
$ws.Range("F18").Value = 1.049166844970305
$ws.Range("I18").Value = 1.039072045662248
$ws.Range("J18").Value = 1.03965908105846
$ws.Range("K18").Value = 1.042219742399668
$ws.Range("L18").Value = 1.044821203887161
$ws.Range("M18").Value = 1.052446123467778
$ws.Range("N18").Value = 1.041135516204126
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033616443765413
$ws.Range("D19").Value = 1.038950406721499
$ws.Range("E19").Value = 1.041586560723577
$ws.Range("F19").Value = 1.049250036158944
$ws.Range("I19").Value = 1.039089922009311
$ws.Range("J19").Value = 1.039707966064779
$ws.Range("K19").Value = 1.042249966787658
$ws.Range("L19").Value = 1.044877177823459
$ws.Range("M19").Value = 1.052514935640489
$ws.Range("N19").Value = 1.041184470632762
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033274796537264
$ws.Range("D20").Value = 1.038750757627182
$ws.Range("E20").Value = 1.041271619388134
$ws.Range("F20").Value = 1.048878057789229
$ws.Range("I20").Value = 1.039009855286529
$ws.Range("J20").Value = 1.03948930700797
$ws.Range("K20").Value = 1.042114716397689
$ws.Range("L20").Value = 1.044626846155054
$ws.Range("M20").Value = 1.052207216371751
$ws.Range("N20").Value = 1.040965501055012
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032164300644106
$ws.Range("D21").Value = 1.038101271966849
$ws.Range("E21").Value = 1.040248544333599
$ws.Range("F21").Value = 1.047669990180872
$ws.Range("I21").Value = 1.038747455227606
$ws.Range("J21").Value = 1.038777834019754
$ws.Range("K21").Value = 1.04167360750464
$ws.Range("L21").Value = 1.043812942490379
$ws.Range("M21").Value = 1.05120723185506
$ws.Range("N21").Value = 1.040253017693545
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031466267554265
$ws.Range("D22").Value = 1.037692610511862
$ws.Range("E22").Value = 1.039605936394714
$ws.Range("F22").Value = 1.046911406754745
$ws.Range("I22").Value = 1.038580870755872
$ws.Range("J22").Value = 1.038330052337474
$ws.Range("K22").Value = 1.041395199647816
$ws.Range("L22").Value = 1.04330117274238
$ws.Range("M22").Value = 1.050578842630781
$ws.Range("N22").Value = 1.039804600109915
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031836290952745
$ws.Range("D23").Value = 1.037909277981706
$ws.Range("E23").Value = 1.039946534102663
$ws.Range("F23").Value = 1.047313453405729
$ws.Range("I23").Value = 1.038669332101279
$ws.Range("J23").Value = 1.038567472387298
$ws.Range("K23").Value = 1.041542889534685
$ws.Range("L23").Value = 1.04357247480431
$ws.Range("M23").Value = 1.050911931154876
$ws.Range("N23").Value = 1.040042357323442
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033293417171787
$ws.Range("D24").Value = 1.038761641071144
$ws.Range("E24").Value = 1.041288782169006
$ws.Range("F24").Value = 1.048898327714723
$ws.Range("I24").Value = 1.039014227238878
$ws.Range("J24").Value = 1.039501227272198
$ws.Range("K24").Value = 1.042122093511386
$ws.Range("L24").Value = 1.044640490701089
$ws.Range("M24").Value = 1.052223986978269
$ws.Range("N24").Value = 1.040977438247382
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034985102527509
$ws.Range("D25").Value = 1.039749384584105
$ws.Range("E25").Value = 1.042849146784635
$ws.Range("F25").Value = 1.050741718210312
$ws.Range("I25").Value = 1.039407458620556
$ws.Range("J25").Value = 1.040582818282619
$ws.Range("K25").Value = 1.042789553583367
$ws.Range("L25").Value = 1.045879686156338
$ws.Range("M25").Value = 1.053748018617454
$ws.Range("N25").Value = 1.042060565241098
